# Module 2 finished, carry on with module 3
# ------------------------------------------------------------------
# This script reproduces the authoring changes described by the
# commit: it reshuffles the "Balance Sheet" summary box, renames a
# handful of bank / asset-class records on "Records" and
# "Accounts & Wealth", adds new per-bank helper headers on "Records",
# and finally leaves "Income" selected as the active sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Balance Sheet
# ------------------------------------------------------------------
$balance = $wb.Worksheets.Item("Balance Sheet")

# "Total Expenses" row label becomes "Cost & Expenses"
$balance.Range("A4").Value = "Cost & Expenses"

# Re-purpose the little scoreboard in P1:P8 so it reports Gross
# Income / Expenses / Net Income instead of just Net Income.
$balance.Range("P1").Value = "Total Gross Income"
$balance.Range("P2").Formula = "=N3"

$balance.Range("P1:P2").Copy()
$balance.Range("P4:P5").PasteSpecial(-4122)
$balance.Range("P4").Value = "Total Expenses"
$balance.Range("P5").Formula = "=N4"

$balance.Range("P1:P2").Copy()
$balance.Range("P7:P8").PasteSpecial(-4122)
$balance.Range("P7").Value = "Total Net Income"
$balance.Range("P8").Formula = "=N7"

# Wealth Allocation class names
$balance.Range("A14").Value = "Class 1"
$balance.Range("A15").Value = "Class 2"
$balance.Range("A16").Value = "Fixed Asset"

# Mirror those class names further down the sheet
$balance.Range("A22").Formula = "=A14"
$balance.Range("A23").Formula = "=A15"
$balance.Range("A24").Formula = "=A16"

# Match the boxed-border look of A23 across A24:A27
$balance.Range("A23").Copy()
$balance.Range("A24:A27").PasteSpecial(-4122)
$balance.Range("A24").Formula = "=A16"

# ------------------------------------------------------------------
# 2) Records - add per-bank helper headers (G1:N1)
# ------------------------------------------------------------------
$records = $wb.Worksheets.Item("Records")

$records.Range("A1").Copy()
$records.Range("G1:N1").PasteSpecial(-4122)

$records.Range("G1").Value = "L1"
$records.Range("H1").Value = "F1"
$records.Range("I1").Value = "L2"
$records.Range("J1").Value = "F4"
$records.Range("K1").Value = "L1"
$records.Range("L1").Value = "F1"
$records.Range("M1").Value = "L2"
$records.Range("N1").Value = "F4"

# ------------------------------------------------------------------
# 3) Accounts & Wealth - renamed banks / classes + refreshed balances
# ------------------------------------------------------------------
$accounts = $wb.Worksheets.Item("Accounts & Wealth")

# Row 2: was "Bank Test 1" / Liquid / 500 / 500 ; Liq1 / 400 / 400 ; Liquid total 750
$accounts.Range("A2").Value = "L1"
$accounts.Range("B2").Value = "Liquid"
$accounts.Range("C2").Value = 2500
$accounts.Range("D2").Value = 2500
$accounts.Range("H2").Value = "Class 1"
$accounts.Range("I2").Value = 2500
$accounts.Range("J2").Value = 2500
$accounts.Range("N2").Value = "Liquid"
$accounts.Range("O2").Value = 3133

# Row 3: was "Bank test 2" / Liquid / 250 / 250 ; Liq2 / 350 / 350 ; Fixed total 2000
$accounts.Range("A3").Value = "F1"
$accounts.Range("B3").Value = "Fixed"
$accounts.Range("C3").Value = 450
$accounts.Range("D3").Value = 450
$accounts.Range("H3").Value = "Class 2"
$accounts.Range("I3").Value = 633
$accounts.Range("J3").Value = 633
$accounts.Range("N3").Value = "Fixed"
$accounts.Range("O3").Value = 461

# Row 4: was "Fixed 1" / Fixed / 2000 / 2000 ; Fixed Asset / 2000 / 2000
$accounts.Range("A4").Value = "L2"
$accounts.Range("B4").Value = "Liquid"
$accounts.Range("C4").Value = 633
$accounts.Range("D4").Value = 633
$accounts.Range("H4").Value = "Fixed Asset"
$accounts.Range("I4").Value = 461
$accounts.Range("J4").Value = 461

# Row 5 is brand new: F4 / Fixed / 11 / 11
$accounts.Range("A4:D4").Copy()
$accounts.Range("A5:D5").PasteSpecial(-4122)
$accounts.Range("A5").Value = "F4"
$accounts.Range("B5").Value = "Fixed"
$accounts.Range("C5").Value = 11
$accounts.Range("D5").Value = 11

# ------------------------------------------------------------------
# 4) Leave "Income" as the active sheet/tab
# ------------------------------------------------------------------
$income = $wb.Worksheets.Item("Income")
$income.Activate()
